$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPDS-JP")

# --- New labels in existing rows 21/23/24/28 ---
# Order matters: it determines the order new entries are appended to the
# shared-string table, so write them in the same sequence as the target file.
$ws.Range("A21").Value = "Abyss Actor - Extra"
$ws.Range("A35").Value = "Casting out the Darklords"
$ws.Range("A37").Value = "Darklords Falling from Grace"
$ws.Range("A23").Value = "Abyss Actor - Wild Hope"
$ws.Range("A24").Value = "Abyss Script - Fantasy Magic"
$ws.Range("A28").Value = "Abyss Script - Rise of the Dark Ruler"

# --- New rows 33-37 (B = card id, C = ":", E = ",") ---
$ws.Range("B33").Value = 100405032
$ws.Range("C33").Value = ":"
$ws.Range("E33").Value = ","

$ws.Range("B34").Value = 100405033
$ws.Range("C34").Value = ":"
$ws.Range("E34").Value = ","

$ws.Range("B35").Value = 100405034
$ws.Range("C35").Value = ":"
$ws.Range("E35").Value = ","

$ws.Range("B36").Value = 100405035
$ws.Range("C36").Value = ":"
$ws.Range("E36").Value = ","

$ws.Range("B37").Value = 100405036
$ws.Range("C37").Value = ":"
$ws.Range("E37").Value = ","

# --- View state: scroll/selection to match the saved workbook state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D36").Select()
